# Update "想去人数" (want-to-go count) figures across sheets, output generated at 456a3b4

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F15").Value = 509
$ws1.Range("F18").Value = 56
$ws1.Range("F19").Value = 433
$ws1.Range("F23").Value = 54804
$ws1.Range("F24").Value = 54804
$ws1.Range("F25").Value = 4346
$ws1.Range("F33").Value = 2890
$ws1.Range("F38").Value = 1156

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 166
$ws2.Range("F16").Value = 7442

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F9").Value = 9317

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F18").Value = 433
$ws4.Range("F21").Value = 54804
$ws4.Range("F22").Value = 166
$ws4.Range("F25").Value = 4346
$ws4.Range("F33").Value = 2890
$ws4.Range("F38").Value = 1156
